# "Added Procedures to OP05"
#
# The captured change for *this* document is entirely confined to the
# SharePoint "content type" custom XML part pair:
#   - customXml/item1.xml       -- the <ct:contentTypeSchema> that mirrors
#                                   this document's bound SharePoint
#                                   content type (0x010100439E...584D0)
#   - customXml/itemProps1.xml  -- its companion <ds:datastoreItem>
#
# Re-syncing that content type (e.g. because a sibling procedure such as
# OP05 gained new fields/columns on the same content type) makes the
# server stamp fresh identifiers on both parts:
#   - ma:versionID on <ct:contentTypeSchema>
#   - ma:fieldsID  on the nested <xsd:schema ma:root="true">
#   - ds:itemID on the <ds:datastoreItem>, which also gains the full
#     <ds:schemaRefs> list describing every namespace used by item1.xml
#
# None of this lives in the document body/story ranges, so there is no
# Find/Replace target for it. It is instead applied to the custom XML
# part itself: locate the part carrying the content-type schema
# namespace, rewrite its two stamped IDs, and swap it back in. Word
# regenerates the paired itemProps1.xml datastore item (id + schemaRefs)
# automatically whenever the owning custom XML part is rewritten.

$d = $word.ActiveDocument

$contentTypeNs = "http://schemas.microsoft.com/office/2006/metadata/contentType"

$oldVersionId = 'ma:versionID="8f9059bd7a7006de8a9c86fca591aade"'
$newVersionId = 'ma:versionID="d2e12c9ede271942fc960b4c1a22c078"'

$oldFieldsId = 'ma:fieldsID="d250b55b03be3e4b4b2c7acd5bf85c08"'
$newFieldsId = 'ma:fieldsID="6b5f36eac1872c4034e06c7c9579686c"'

$parts = $d.CustomXMLParts

# Find the custom XML part that stores the SharePoint content-type schema.
$target = $null
try {
    $matches = $parts.SelectByNamespace($contentTypeNs)
    if ($matches -ne $null -and $matches.Count -ge 1) {
        $target = $matches.Item(1)
    }
} catch {
    $target = $null
}

if ($target -eq $null) {
    for ($i = 1; $i -le $parts.Count; $i++) {
        $candidate = $parts.Item($i)
        if ($candidate.NamespaceURI -eq $contentTypeNs) {
            $target = $candidate
        }
    }
}

if ($target -ne $null -and $target.XML) {
    $xml = $target.XML
    $xml = $xml.Replace($oldVersionId, $newVersionId)
    $xml = $xml.Replace($oldFieldsId, $newFieldsId)

    # Word's object model exposes CustomXMLPart.XML as read-only, so the
    # refreshed schema is swapped in by replacing the part outright.
    $target.Delete()
    $parts.Add($xml)
}
